$d = $word.ActiveDocument

# The bibliography list ends with a paragraph that has a hyperlink to the
# nature.com s41598 article, immediately followed by one last empty list
# paragraph. We need to insert three brand-new list paragraphs between
# them.
$lastHyperlinkPara = $d.Paragraphs(26)
$anchor = $d.Range($lastHyperlinkPara.Range.End, $lastHyperlinkPara.Range.End)

# Creating three paragraph breaks at the same anchor point inserts three
# new (empty) paragraphs right after paragraph 26, pushing the old final
# empty paragraph further down - exactly where we want our new entries.
$anchor.InsertParagraphAfter()
$anchor.InsertParagraphAfter()
$anchor.InsertParagraphAfter()

# --- New paragraph 1: hyperlink to the Science Mag article ---
$p1 = $d.Paragraphs(27)
$p1Start = $d.Range($p1.Range.Start, $p1.Range.Start)
$url1 = "https://www.sciencemag.org/news/2015/10/sensors-may-soon-give-prosthetics-lifelike-sense-touch"
$d.Hyperlinks.Add($p1Start, $url1, "", "", $url1) | Out-Null
$p1 = $d.Paragraphs(27)
$p1InsertPoint = $p1.Range.End - 1
$d.Range($p1InsertPoint, $p1InsertPoint).InsertAfter(" ")

# --- New paragraph 2: hyperlink to the Rehabmart product page ---
$p2 = $d.Paragraphs(28)
$p2Start = $d.Range($p2.Range.Start, $p2.Range.Start)
$url2 = "https://www.rehabmart.com/product/smart-glove-for-stroke-rehabilitation-by-neofect-49247.html?gclid=Cj0KCQiA3NX_BRDQARIsALA3fILmFJC7fMqFPAU3qoxPEqXB0Ly_DZAwrzM9IlDR2tFiQVmghzk41lYaAuGyEALw_wcB"
$d.Hyperlinks.Add($p2Start, $url2, "", "", $url2) | Out-Null
$p2 = $d.Paragraphs(28)
$p2InsertPoint = $p2.Range.End - 1
$d.Range($p2InsertPoint, $p2InsertPoint).InsertAfter(" ")

# --- New paragraph 3: plain DOI citation text ---
$p3 = $d.Paragraphs(29)
$p3Start = $d.Range($p3.Range.Start, $p3.Range.Start)
$p3Start.InsertAfter("DOI: 10.1016/j.mejo.2018.01.014")
$p3 = $d.Paragraphs(29)
$p3InsertPoint = $p3.Range.End - 1
$leftDoubleQuote = [char]0x201C
$d.Range($p3InsertPoint, $p3InsertPoint).InsertAfter(" ${leftDoubleQuote}Wearable technologies for hand joints monitoring for rehabilitation")

Write-Host "Done"
